# TC04_Canine_Filter_Diagnosis-Melanoma.xlsx
#
# The "CasesTab" row's Cypher query (cell B2 on the single "startup" sheet)
# included a `Cohort` column (joined through a `cohort` node) that is being
# dropped from the report. Remove the trailing
#     coalesce(co.cohort_description, '') AS `Cohort`
# line (and the now-superfluous trailing comma on the previous line) from
# that query text, leaving the query ending on the "Response to Treatment"
# column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (co:cohort)<-[*]-(c)
WHERE diag.disease_term IN ['Melanoma']
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value2 = $newCasesQuery

# The file was also re-saved with the view scrolled back to the top and the
# selection sitting on B2 (rather than scrolled down to B4) with the zoom
# reset to 100%.
[void]$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100

Write-Host "B2 is now:"
Write-Host $ws.Range("B2").Value2
